$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6975.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 6975.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 6975.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -7113.5
$ws.Range("H51").Value = 2049.6155
$ws.Range("J51").Value = 2244.5454
$ws.Range("L51").Value = 2244.5454
$ws.Range("N51").Value = -3212.5454
$ws.Range("H96").Value = 651.4706
$ws.Range("I96").Value = 389
$ws.Range("J96").Value = 1504.5
$ws.Range("K96").Value = 1167
$ws.Range("L96").Value = 4513.5
$ws.Range("M96").Value = 206
$ws.Range("N96").Value = -7259.5
$ws.Range("H101").Value = 671
$ws.Range("I101").Value = 330.77777
$ws.Range("K101").Value = 992.33331
$ws.Range("M101").Value = 629.66669
$ws.Range("H106").Value = 1030.091
$ws.Range("I106").Value = 1040.625
$ws.Range("J106").Value = 1002
$ws.Range("K106").Value = 1040.625
$ws.Range("L106").Value = 1002
$ws.Range("M106").Value = -409.625
$ws.Range("N106").Value = -2264
$ws.Range("H107").Value = 587.8570999999999
$ws.Range("I107").Value = 546.7646999999999
$ws.Range("K107").Value = 546.7646999999999
$ws.Range("M107").Value = 1373.2353
$ws.Range("H129").Value = 17860988
$ws.Range("I129").Value = 41667948
$ws.Range("J129").Value = 5766.75
$ws.Range("K129").Value = 125003844
$ws.Range("L129").Value = 17300.25
$ws.Range("M129").Value = -124998844
$ws.Range("N129").Value = -27300.25
$ws.Range("H137").Value = 5888100.5
$ws.Range("I137").Value = 14295171
$ws.Range("J137").Value = 3150.8
$ws.Range("K137").Value = 42885513
$ws.Range("L137").Value = 9452.400000000001
$ws.Range("M137").Value = -42882963
$ws.Range("N137").Value = -14552.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8330.681
$ws.Range("I32").Value = 5539.6665
$ws.Range("J32").Value = 21936.875
$ws.Range("K32").Value = 5539.6665
$ws.Range("L32").Value = 21936.875
$ws.Range("M32").Value = -5252.6665
$ws.Range("N32").Value = -22510.875
$ws.Range("H45").Value = 1687.0938
$ws.Range("I45").Value = 1067.1666
$ws.Range("J45").Value = 3546.875
$ws.Range("K45").Value = 1067.1666
$ws.Range("L45").Value = 3546.875
$ws.Range("M45").Value = -690.1666
$ws.Range("N45").Value = -4300.875
$ws.Range("H74").Value = 827.8182
$ws.Range("I74").Value = 681.2778
$ws.Range("J74").Value = 1487.25
$ws.Range("K74").Value = 681.2778
$ws.Range("L74").Value = 1487.25
$ws.Range("M74").Value = 192.7222
$ws.Range("N74").Value = -3235.25
$ws.Range("H77").Value = 827.8182
$ws.Range("I77").Value = 681.2778
$ws.Range("J77").Value = 1487.25
$ws.Range("K77").Value = 3406.389
$ws.Range("L77").Value = 7436.25
$ws.Range("M77").Value = 961.6110000000003
$ws.Range("N77").Value = -16172.25
$ws.Range("H97").Value = 418.36
$ws.Range("I97").Value = 373.45
$ws.Range("J97").Value = 598
$ws.Range("K97").Value = 373.45
$ws.Range("L97").Value = 598
$ws.Range("M97").Value = 122.55
$ws.Range("N97").Value = -1590
$ws.Range("H102").Value = 2714.8
$ws.Range("I102").Value = 2320.4348
$ws.Range("K102").Value = 2320.4348
$ws.Range("M102").Value = -698.4348
$ws.Range("H122").Value = 3496
$ws.Range("J122").Value = 3496
$ws.Range("L122").Value = 10488
$ws.Range("N122").Value = -15388

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1724.92
$ws.Range("I20").Value = 1168.5
$ws.Range("J20").Value = 2433.0908
$ws.Range("K20").Value = 1168.5
$ws.Range("L20").Value = 2433.0908
$ws.Range("M20").Value = -921.5
$ws.Range("N20").Value = -2927.0908
$ws.Range("H80").Value = 908.46155
$ws.Range("I80").Value = 1010.36365
$ws.Range("J80").Value = 833.73334
$ws.Range("K80").Value = 1010.36365
$ws.Range("L80").Value = 833.73334
$ws.Range("M80").Value = -12.36365000000001
$ws.Range("N80").Value = -2829.73334
$ws.Range("H83").Value = 908.46155
$ws.Range("I83").Value = 1010.36365
$ws.Range("J83").Value = 833.73334
$ws.Range("K83").Value = 5051.81825
$ws.Range("L83").Value = 4168.6667
$ws.Range("M83").Value = -59.81825000000026
$ws.Range("N83").Value = -14152.6667
$ws.Range("H86").Value = 2122.0435
$ws.Range("I86").Value = 1270
$ws.Range("J86").Value = 2777.4614
$ws.Range("K86").Value = 1270
$ws.Range("L86").Value = 2777.4614
$ws.Range("M86").Value = -147
$ws.Range("N86").Value = -5023.4614
$ws.Range("H89").Value = 2122.0435
$ws.Range("I89").Value = 1270
$ws.Range("J89").Value = 2777.4614
$ws.Range("K89").Value = 6350
$ws.Range("L89").Value = 13887.307
$ws.Range("M89").Value = -734
$ws.Range("N89").Value = -25119.307
$ws.Range("H99").Value = 1626.1666
$ws.Range("I99").Value = 915.619
$ws.Range("J99").Value = 6600
$ws.Range("K99").Value = 915.619
$ws.Range("L99").Value = 6600
$ws.Range("M99").Value = 582.381
$ws.Range("N99").Value = -9596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2086300.1
$ws.Range("I31").Value = 2565472
$ws.Range("J31").Value = 9888.888999999999
$ws.Range("K31").Value = 2565472
$ws.Range("L31").Value = 9888.888999999999
$ws.Range("M31").Value = -2565177
$ws.Range("N31").Value = -10478.889
$ws.Range("H34").Value = 2086300.1
$ws.Range("I34").Value = 2565472
$ws.Range("J34").Value = 9888.888999999999
$ws.Range("K34").Value = 2565472
$ws.Range("L34").Value = 9888.888999999999
$ws.Range("M34").Value = -2565270
$ws.Range("N34").Value = -10292.889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 23000
$ws.Range("J74").Value = 28000
$ws.Range("L74").Value = 84000
$ws.Range("N74").Value = -86122
$ws.Range("H77").Value = 23000
$ws.Range("J77").Value = 28000
$ws.Range("L77").Value = 252000
$ws.Range("N77").Value = -262608
$ws.Range("H137").Value = 3287.0715
$ws.Range("I137").Value = 2754.75
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 8264.25
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -3164.25
$ws.Range("N137").Value = -20700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1979.619
$ws.Range("I102").Value = 1660.9231
$ws.Range("J102").Value = 2497.5
$ws.Range("K102").Value = 1660.9231
$ws.Range("L102").Value = 2497.5
$ws.Range("M102").Value = -38.92309999999998
$ws.Range("N102").Value = -5741.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2282
$ws.Range("I40").Value = 1130.4
$ws.Range("J40").Value = 3433.6
$ws.Range("K40").Value = 1130.4
$ws.Range("L40").Value = 3433.6
$ws.Range("M40").Value = -994.4000000000001
$ws.Range("N40").Value = -3705.6
$ws.Range("H93").Value = 1890.08
$ws.Range("I93").Value = 1241.5555
$ws.Range("J93").Value = 3557.7144
$ws.Range("K93").Value = 1241.5555
$ws.Range("L93").Value = 3557.7144
$ws.Range("M93").Value = 6.444500000000062
$ws.Range("N93").Value = -6053.7144
$ws.Range("H129").Value = 26000
$ws.Range("J129").Value = 26000
$ws.Range("L129").Value = 26000
$ws.Range("N129").Value = -36000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1514.1428
$ws.Range("I96").Value = 1450
$ws.Range("J96").Value = 1539.8
$ws.Range("K96").Value = 1450
$ws.Range("L96").Value = 1539.8
$ws.Range("M96").Value = -77
$ws.Range("N96").Value = -4285.8
$ws.Range("H100").Value = 462.2857
$ws.Range("I100").Value = 420.92307
$ws.Range("K100").Value = 841.84614
$ws.Range("M100").Value = -300.84614
